$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for the "Resolving-Mac" sending-cluster block (old rows 14-17),
# which no longer exist after the TPM update.
$ws.Range("A14:T17").Delete()

$data = @(
  @("ECs", "Sema3a", "Nrp2", "ECs", 3, 1, 1.339548666666667, 4.018646, 0.09827984122213275, 0.09827984122213274, 3, 1, 70.23436, 210.70308, 0.7023186840741513, 0.7023186840741513, 94.08234329218666, 846.7410896296801, 0.06902376875814481, 0.0690237687581448),
  @("ECs", "Sema3a", "Nrp2", "FAPs", 3, 1, 1.339548666666667, 4.018646, 0.09827984122213275, 0.09827984122213274, 3, 1, 6.551362, 19.654086, 0.06551129587759326, 0.06551129587759325, 8.775868231950668, 78.98281408755601, 0.006438439757106025, 0.006438439757106023),
  @("ECs", "Sema3a", "Nrp2", "MuSCs", 3, 1, 1.339548666666667, 4.018646, 0.09827984122213275, 0.09827984122213274, 3, 1, 4.268944666666666, 12.806834, 0.04268793224112385, 0.04268793224112385, 5.718459136307111, 51.466132226764, 0.004195363202758814, 0.004195363202758813),
  @("ECs", "Sema3a", "Nrp2", "Resolving-Mac", 3, 1, 1.339548666666667, 4.018646, 0.09827984122213275, 0.09827984122213274, 3, 1, 18.948881, 56.846643, 0.1894820878071316, 0.1894820878071315, 25.38294827837533, 228.446534505378, 0.01862226950412311, 0.0186222695041231),
  @("FAPs", "Sema3a", "Nrp2", "ECs", 3, 1, 2.040291, 6.120873, 0.1496918182345096, 0.1496918182345096, 3, 1, 70.23436, 210.70308, 0.7023186840741513, 0.7023186840741513, 143.29853259876, 1289.68679338884, 0.1051313607991279, 0.1051313607991279),
  @("FAPs", "Sema3a", "Nrp2", "FAPs", 3, 1, 2.040291, 6.120873, 0.1496918182345096, 0.1496918182345096, 3, 1, 6.551362, 19.654086, 0.06551129587759326, 0.06551129587759325, 13.366684926342, 120.300164337078, 0.009806504994815871, 0.00980650499481587),
  @("FAPs", "Sema3a", "Nrp2", "MuSCs", 3, 1, 2.040291, 6.120873, 0.1496918182345096, 0.1496918182345096, 3, 1, 4.268944666666666, 12.806834, 0.04268793224112385, 0.04268793224112385, 8.709889382897998, 78.38900444608198, 0.006390034193845376, 0.006390034193845376),
  @("FAPs", "Sema3a", "Nrp2", "Resolving-Mac", 3, 1, 2.040291, 6.120873, 0.1496918182345096, 0.1496918182345096, 3, 1, 18.948881, 56.846643, 0.1894820878071316, 0.1894820878071315, 38.661231364371, 347.951082279339, 0.02836391824672053, 0.02836391824672053),
  @("MuSCs", "Sema3a", "Nrp2", "ECs", 3, 1, 10.25010366666667, 30.750311, 0.7520283405433575, 0.7520283405433575, 3, 1, 70.23436, 210.70308, 0.7023186840741513, 0.7023186840741513, 719.9094709619866, 6479.18523865788, 0.5281635545168786, 0.5281635545168786),
  @("MuSCs", "Sema3a", "Nrp2", "FAPs", 3, 1, 10.25010366666667, 30.750311, 0.7520283405433575, 0.7520283405433575, 3, 1, 6.551362, 19.654086, 0.06551129587759326, 0.06551129587759325, 67.15213965786066, 604.369256920746, 0.04926635112567136, 0.04926635112567135),
  @("MuSCs", "Sema3a", "Nrp2", "MuSCs", 3, 1, 10.25010366666667, 30.750311, 0.7520283405433575, 0.7520283405433575, 3, 1, 4.268944666666666, 12.806834, 0.04268793224112385, 0.04268793224112385, 43.7571253805971, 393.814128425374, 0.03210253484451966, 0.03210253484451966),
  @("MuSCs", "Sema3a", "Nrp2", "Resolving-Mac", 3, 1, 10.25010366666667, 30.750311, 0.7520283405433575, 0.7520283405433575, 3, 1, 18.948881, 56.846643, 0.1894820878071316, 0.1894820878071315, 194.2279946173303, 1748.051951555973, 0.1424959000562879, 0.1424959000562879),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}
